$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" - record a new handoff timestamp for the
# file 7d78cd29-f389-403f-ab08-5219f03e017f (row 7 on each sheet).

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-29-20 18:29:34"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-20 18:29:31"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-20 18:29:34"
